$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7886165380477905
$ws.Range("B1").Value = 1.001588940620422
$ws.Range("C1").Value = 1.478047132492065
$ws.Range("D1").Value = 2.203947067260742
$ws.Range("E1").Value = 1.625526785850525
